$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 263.83334
$ws.Range("I33").Value = 259.1875
$ws.Range("K33").Value = 259.1875
$ws.Range("M33").Value = -30.1875
$ws.Range("H129").Value = 1055.0968
$ws.Range("J129").Value = 1113.6296
$ws.Range("L129").Value = 3340.8888
$ws.Range("N129").Value = -13340.8888
$ws.Range("H133").Value = 69770
$ws.Range("J133").Value = 69770
$ws.Range("L133").Value = 69770
$ws.Range("N133").Value = -79890
$ws.Range("H136").Value = 33118.75
$ws.Range("J136").Value = 33118.75
$ws.Range("L136").Value = 33118.75
$ws.Range("N136").Value = -43318.75
$ws.Range("H138").Value = 2266.32
$ws.Range("J138").Value = 2276.908
$ws.Range("L138").Value = 6830.724
$ws.Range("N138").Value = -17110.724
$ws.Range("H139").Value = 66125
$ws.Range("J139").Value = 66125
$ws.Range("L139").Value = 66125
$ws.Range("N139").Value = -76405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 6039960
$ws.Range("I13").Value = 15025000
$ws.Range("K13").Value = 15025000
$ws.Range("M13").Value = -15024856
$ws.Range("H32").Value = 12016.577
$ws.Range("I32").Value = 11866.318
$ws.Range("K32").Value = 11866.318
$ws.Range("M32").Value = -11579.318
$ws.Range("H103").Value = 177500
$ws.Range("J103").Value = 177500
$ws.Range("L103").Value = 177500
$ws.Range("N103").Value = -179844
$ws.Range("H112").Value = 23699.4
$ws.Range("J112").Value = 23699.4
$ws.Range("L112").Value = 23699.4
$ws.Range("N112").Value = -26653.4
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H135").Value = 41500
$ws.Range("J135").Value = 41500
$ws.Range("L135").Value = 41500
$ws.Range("N135").Value = -51640

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 47400
$ws.Range("J59").Value = 47400
$ws.Range("L59").Value = 47400
$ws.Range("N59").Value = -49094
$ws.Range("H109").Value = 26156.111
$ws.Range("J109").Value = 26156.111
$ws.Range("L109").Value = 26156.111
$ws.Range("N109").Value = -28930.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2395
$ws.Range("I5").Value = 465
$ws.Range("K5").Value = 465
$ws.Range("M5").Value = -353
$ws.Range("H22").Value = 265.57144
$ws.Range("J22").Value = 460.5
$ws.Range("L22").Value = 460.5
$ws.Range("N22").Value = -1160.5
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 19610442
$ws.Range("I31").Value = 29413554
$ws.Range("J31").Value = 4217.647
$ws.Range("K31").Value = 29413554
$ws.Range("L31").Value = 4217.647
$ws.Range("M31").Value = -29413259
$ws.Range("N31").Value = -4807.647
$ws.Range("H34").Value = 19610442
$ws.Range("I34").Value = 29413554
$ws.Range("J34").Value = 4217.647
$ws.Range("K34").Value = 29413554
$ws.Range("L34").Value = 4217.647
$ws.Range("M34").Value = -29413352
$ws.Range("N34").Value = -4621.647
$ws.Range("H39").Value = 24550.334
$ws.Range("I39").Value = 24550.334
$ws.Range("K39").Value = 24550.334
$ws.Range("M39").Value = -24159.334
$ws.Range("H49").Value = 24550.334
$ws.Range("I49").Value = 24550.334
$ws.Range("K49").Value = 24550.334
$ws.Range("M49").Value = -24368.334
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H99").Value = 4375.8335
$ws.Range("I99").Value = 4451
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 4451
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -2953
$ws.Range("N99").Value = -6996
$ws.Range("H126").Value = 4375.8335
$ws.Range("I126").Value = 4451
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 13353
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -10883
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6239.125
$ws.Range("I75").Value = 2637.6667
$ws.Range("J75").Value = 8400
$ws.Range("K75").Value = 7913.000100000001
$ws.Range("L75").Value = 25200
$ws.Range("M75").Value = -6915.000100000001
$ws.Range("N75").Value = -27196
$ws.Range("H78").Value = 6239.125
$ws.Range("I78").Value = 2637.6667
$ws.Range("J78").Value = 8400
$ws.Range("K78").Value = 23739.0003
$ws.Range("L78").Value = 75600
$ws.Range("M78").Value = -18747.0003
$ws.Range("N78").Value = -85584
$ws.Range("H113").Value = 875
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 914.2857
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2742.8571
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -7082.8571
$ws.Range("H122").Value = 758.03705
$ws.Range("I122").Value = 498.21054
$ws.Range("J122").Value = 1375.125
$ws.Range("K122").Value = 4483.894859999999
$ws.Range("L122").Value = 12376.125
$ws.Range("M122").Value = -2033.894859999999
$ws.Range("N122").Value = -17276.125
$ws.Range("H125").Value = 2846.4285
$ws.Range("I125").Value = 1818.3334
$ws.Range("J125").Value = 3617.5
$ws.Range("K125").Value = 5455.0002
$ws.Range("L125").Value = 10852.5
$ws.Range("M125").Value = -535.0002000000004
$ws.Range("N125").Value = -20692.5
$ws.Range("H130").Value = 2450
$ws.Range("I130").Value = 1000
$ws.Range("J130").Value = 2933.3333
$ws.Range("K130").Value = 3000
$ws.Range("L130").Value = 8799.999899999999
$ws.Range("M130").Value = 2020
$ws.Range("N130").Value = -18839.9999
$ws.Range("H133").Value = 4253.5713
$ws.Range("J133").Value = 5973.3335
$ws.Range("L133").Value = 17920.0005
$ws.Range("N133").Value = -28040.0005
$ws.Range("H134").Value = 7316.9067
$ws.Range("I134").Value = 8924
$ws.Range("J134").Value = 6455.9644
$ws.Range("K134").Value = 26772
$ws.Range("L134").Value = 19367.8932
$ws.Range("M134").Value = -21702
$ws.Range("N134").Value = -29507.8932
$ws.Range("H137").Value = 47621320
$ws.Range("J137").Value = 83336590
$ws.Range("L137").Value = 250009770
$ws.Range("N137").Value = -250019970
$ws.Range("H139").Value = 1927.2333
$ws.Range("I139").Value = 1511.409
$ws.Range("J139").Value = 3070.75
$ws.Range("K139").Value = 4534.227000000001
$ws.Range("L139").Value = 9212.25
$ws.Range("M139").Value = 605.7729999999992
$ws.Range("N139").Value = -19492.25
$ws.Range("H141").Value = 5564.7144
$ws.Range("I141").Value = 5820
$ws.Range("K141").Value = 17460
$ws.Range("M141").Value = -12280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2931.6924
$ws.Range("I102").Value = 3011.2
$ws.Range("J102").Value = 2666.6667
$ws.Range("K102").Value = 3011.2
$ws.Range("L102").Value = 2666.6667
$ws.Range("M102").Value = -1389.2
$ws.Range("N102").Value = -5910.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 3500
$ws.Range("I35").Value = 3500
$ws.Range("K35").Value = 3500
$ws.Range("M35").Value = -3164
$ws.Range("H122").Value = 4976.2354
$ws.Range("I122").Value = 4828.2856
$ws.Range("J122").Value = 5666.6665
$ws.Range("K122").Value = 14484.8568
$ws.Range("L122").Value = 16999.9995
$ws.Range("M122").Value = -12034.8568
$ws.Range("N122").Value = -21899.9995
$ws.Range("H135").Value = 192462.25
$ws.Range("J135").Value = 192462.25
$ws.Range("L135").Value = 192462.25
$ws.Range("N135").Value = -202602.25
$ws.Range("H138").Value = 25000
$ws.Range("J138").Value = 25000
$ws.Range("L138").Value = 25000
$ws.Range("N138").Value = -35280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 9859710
$ws.Range("J39").Value = 14966.667
$ws.Range("L39").Value = 14966.667
$ws.Range("N39").Value = -15792.667
$ws.Range("H43").Value = 49999
$ws.Range("I43").Value = 49999
$ws.Range("K43").Value = 49999
$ws.Range("M43").Value = -49850
$ws.Range("H81").Value = 45058.305
$ws.Range("I81").Value = 60038.883
$ws.Range("K81").Value = 120077.766
$ws.Range("M81").Value = -119016.766
$ws.Range("H84").Value = 45058.305
$ws.Range("I84").Value = 60038.883
$ws.Range("K84").Value = 600388.8300000001
$ws.Range("M84").Value = -595084.8300000001
$ws.Range("H137").Value = 34993.332
$ws.Range("J137").Value = 34993.332
$ws.Range("L137").Value = 34993.332
$ws.Range("N137").Value = -45193.332
$ws.Range("H139").Value = 62302.223
$ws.Range("J139").Value = 62302.223
$ws.Range("L139").Value = 62302.223
$ws.Range("N139").Value = -72582.223
